$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parameters")
$ws.Range("B6").Value = 1407
$ws.Range("C6").Value = 62
$ws.Range("C6").NumberFormat = "0.00E+00"
